$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "GenCor2019"
$ws.Range("B2").Value = "Elecciones Generales de la Provincia de Córdoba"
$ws.Range("C2").Value = "Gobernador y Vice Gobernador"

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "2019-05-12"
$ws.Range("E2").Value = "2019-05-12"
$ws.Range("D2:E2").ClearFormats()
